$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)

$para1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>Eye Irritation/Corrosion</w:t>
  </w:r>
</w:p>
'@

$endRange.InsertXML($para1Xml)

$endRange2 = $d.Content
$endRange2.Collapse(0)

$para2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>118 compounds tested for eye irritation/corrosion by Yamaguchi et al.</w:t>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:instrText xml:space="preserve"> ADDIN ZOTERO_ITEM CSL_CITATION {"citationID":"fefntu31","properties":{"formattedCitation":"\\super 3\\nosupersub{}","plainCitation":"3","noteIndex":0},"citationItems":[{"id":3445,"uris":["http://zotero.org/users/11717528/items/AZ6YE7TW"],"itemData":{"id":3445,"type":"article-journal","abstract":"We recently developed a novel Vitrigel-eye irritancy test (EIT) method. The Vitrigel-EIT method is composed of two parts, i.e., the construction of a human corneal epithelium (HCE) model in a collagen vitrigel membrane chamber and the prediction of eye irritancy by analyzing the time-dependent profile of transepithelial electrical resistance values for 3 min after exposing a chemical to the HCE model. In this study, we estimated the predictive performance of Vitrigel-EIT method by testing a total of 118 chemicals. The category determined by the Vitrigel-EIT method in comparison to the globally harmonized system classification revealed that the sensitivity, specificity and accuracy were 90.1%, 65.9% and 80.5%, respectively. Here, five of seven false-negative chemicals were acidic chemicals inducing the irregular rising of transepithelial electrical resistance values. In case of eliminating the test chemical solutions showing pH 5 or lower, the sensitivity, specificity and accuracy were improved to 96.8%, 67.4% and 84.4%, respectively. Meanwhile, nine of 16 false-positive chemicals were classified irritant by the US Environmental Protection Agency. In addition, the disappearance of ZO-1, a tight junction-associated protein and MUC1, a cell membrane-spanning mucin was immunohistologically confirmed in the HCE models after exposing not only eye irritant chemicals but also false-positive chemicals, suggesting that such false-positive chemicals have an eye irritant potential. These data demonstrated that the Vitrigel-EIT method could provide excellent predictive performance to judge the widespread eye irritancy, including very mild irritant chemicals. We hope that the Vitrigel-EIT method contributes to the development of safe commodity chemicals. Copyright © 2015 The Authors. Journal of Applied Toxicology published by John Wiley &amp; Sons Ltd.","container-title":"Journal of Applied Toxicology","DOI":"10.1002/jat.3254","ISSN":"1099-1263","issue":"8","language":"en","license":"Copyright © 2015 The Authors. Journal of Applied Toxicology published by John Wiley &amp; Sons Ltd.","note":"_eprint: https://onlinelibrary.wiley.com/doi/pdf/10.1002/jat.3254","page":"1025-1037","source":"Wiley Online Library","title":"Predictive performance of the Vitrigel-eye irritancy test method using 118 chemicals","volume":"36","author":[{"family":"Yamaguchi","given":"Hiroyuki"},{"family":"Kojima","given":"Hajime"},{"family":"Takezawa","given":"Toshiaki"}],"issued":{"date-parts":[["2016"]]}}}],"schema":"https://github.com/citation-style-language/schema/raw/master/csl-citation.json"} </w:instrText>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>3</w:t>
  </w:r>
  <w:r>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
'@

$endRange2.InsertXML($para2Xml)

Write-Host "Paragraphs count after edit:" $d.Paragraphs.Count
